$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260507583618164
$ws.Range("B1").Value = 3.221070051193237
$ws.Range("C1").Value = 6.003250122070312
$ws.Range("D1").Value = 1.784162282943726
$ws.Range("E1").Value = 1.048321843147278
